# "Add files via upload" — the uploaded workbook now has one more day of
# data appended as row 30 (continuing the daily series from row 29),
# which used to be a blank placeholder row. The view/selection also moved
# down to reflect the newly-added data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new day's figures that used to be an empty row.
$ws.Cells.Item(30, 1).Value = 45979   # date serial for the new day
$ws.Cells.Item(30, 2).Value = 943
$ws.Cells.Item(30, 3).Value = 11
$ws.Cells.Item(30, 4).Value = 932

# Bring the new row into view / reselect where the user was last working.
[void]$ws.Activate()
$win = $excel.Windows.Item(1)
$win.ScrollRow = 20
$win.ScrollColumn = 1
[void]$ws.Range("H33").Select()
